# Add two new key-item rows (CA.23 / CA.25) right after the "CA.5 - GROUP 2"
# row, and close up the stale gap that used to sit between row 26 and row 37
# (the "CA.20 - Fees/ Total asset" row), so the sheet ends up as one
# contiguous block A1:B29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the unused/empty rows 27-36 so the last data row ("Fees/ Total
#    asset", currently row 37) slides up to become row 27 with no gap.
$ws.Rows("27:36").Delete()

# 2) Insert two fresh rows right after row 10 ("CA.5 | GROUP 2") to make
#    room for the new metrics.
$ws.Rows("11:12").Insert()

# 3) Populate the two new rows. Column A first for both rows, then column B,
#    so new shared-strings are appended in KeyCode/Name pairs.
$ws.Range("A11").Value = "CA.23"
$ws.Range("A12").Value = "CA.25"
$ws.Range("B11").Value = "NPL Formation (%)"
$ws.Range("B12").Value = "G2 Formation (%)"

# 4) Match the author's final selection/cursor position.
$ws.Range("B13").Select()
